$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1.78941,1.78941,1.78941,1.78941,1.32543,1.32543,1.78941,1.78941,1.78941,1.78941,
    1.32543,1.32543,1.78941,1.78941,1.78941,1.32543,1.32543,1.78941,1.78941,1.78941,
    1.78941,1.32543,1.32543,1.78941,1.78941,1.78941,1.32543,1.32543,1.78941,1.78941,
    1.78941,1.32543,1.32543,1.32543,1.78941,1.78941,1.78941,1.78941,1.78941,1.32543,
    1.78941,1.78941,1.78941,1.78941,1.32543,1.78941,1.78941,1.78941,1.78941,1.32543,
    1.32543,1.78941,1.78941,1.78941,1.78941,1.32543,1.32543,1.78941,1.78941,1.32543,
    1.32543,1.32543,1.78941,1.78941,1.32543,1.32543,1.32543
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
